{"js": "// Apply the French-translation edits described by the diff.\n// Each entry is an exact, unique source substring found in a run's\n// text, paired with its replacement. We use Body.search (matchCase,\n// whole string) to locate the run's text and Range.insertText(...,\n// \"Replace\") to swap it in-place, preserving run formatting.\nconst replacements = [\n  [\"The administrations of three\", \"Les administrations de trois\"],\n  [\"neighboring cities: A, B and C decided\", \"villes voisines: A, B et C ont d\u00e9cid\u00e9\"],\n  [\"to build an airport dividing the costs of\", \"de construire un a\u00e9roport et de diviser les co\u00fbts de ce\"],\n  [\"00:00:46,000 --> 00:00:48,000\", \"0:00:46,000 --> 00:00:48,000\"],\n  [\"implementation. The condition on the\", \"projet. L'endroit, qui\"],\n  [\"choice of the most suitable place is\", \"convient le plus est\"],\n];\n\nconst body = context.document.body;\n\nfor (const [find, replace] of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + find);\n  }\n\n  results.items[0].insertText(replace, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Apply the French-translation edits described by the diff.\n# For each (find, replace) pair: locate the unique source run text with\n# Range.Find, then assign the replacement directly to the found Range's\n# .Text. Assigning .Text (rather than driving Find.Execute's built-in\n# Replace) keeps the literal characters (e.g. a straight apostrophe)\n# intact instead of letting Word's smart-quote autocorrect rewrite them,\n# and it preserves the run's existing character formatting.\n\n$d = $word.ActiveDocument\n\nfunction Replace-FirstMatch($doc, $findText, $replaceText) {\n    $r = $doc.Content\n    $find = $r.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $found = $find.Execute()\n    if (-not $found) {\n        throw \"Not found: $findText\"\n    }\n    $r.Text = $replaceText\n}\n\nReplace-FirstMatch $d \"The administrations of three\" \"Les administrations de trois\"\nReplace-FirstMatch $d \"neighboring cities: A, B and C decided\" \"villes voisines: A, B et C ont d\u00e9cid\u00e9\"\nReplace-FirstMatch $d \"to build an airport dividing the costs of\" \"de construire un a\u00e9roport et de diviser les co\u00fbts de ce\"\nReplace-FirstMatch $d \"00:00:46,000 --> 00:00:48,000\" \"0:00:46,000 --> 00:00:48,000\"\nReplace-FirstMatch $d \"implementation. The condition on the\" \"projet. L'endroit, qui\"\nReplace-FirstMatch $d \"choice of the most suitable place is\" \"convient le plus est\"\n"}
